$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 570
$ws.Range("C2").Value = 529
$ws.Range("D2").Value = 564
$ws.Range("E2").Value = 567
$ws.Range("F2").Value = 531
$ws.Range("G2").Value = 528
$ws.Range("I2").Value = 854
$ws.Range("J2").Value = 875
$ws.Range("K2").Value = 825
$ws.Range("L2").Value = 793
$ws.Range("M2").Value = 789
$ws.Range("B3").Value = 514
$ws.Range("C3").Value = 537
$ws.Range("D3").Value = 558
$ws.Range("E3").Value = 573
$ws.Range("F3").Value = 583
$ws.Range("G3").Value = 555
$ws.Range("I3").Value = 831
$ws.Range("J3").Value = 891
$ws.Range("K3").Value = 836
$ws.Range("L3").Value = 848
$ws.Range("M3").Value = 763
$ws.Range("B4").Value = 477
$ws.Range("C4").Value = 480
$ws.Range("D4").Value = 471
$ws.Range("E4").Value = 481
$ws.Range("F4").Value = 490
$ws.Range("G4").Value = 481
$ws.Range("H4").Value = 667
$ws.Range("I4").Value = 679
$ws.Range("K4").Value = 658
$ws.Range("L4").Value = 646
$ws.Range("M4").Value = 623
$ws.Range("B5").Value = 425
$ws.Range("C5").Value = 373
$ws.Range("D5").Value = 435
$ws.Range("E5").Value = 419
$ws.Range("F5").Value = 447
$ws.Range("G5").Value = 425
$ws.Range("H5").Value = 625
$ws.Range("I5").Value = 632
$ws.Range("J5").Value = 647
$ws.Range("K5").Value = 613
$ws.Range("M5").Value = 578
$ws.Range("B6").Value = 474
$ws.Range("D6").Value = 445
$ws.Range("E6").Value = 452
$ws.Range("F6").Value = 459
$ws.Range("G6").Value = 402
$ws.Range("H6").Value = 656
$ws.Range("I6").Value = 662
$ws.Range("J6").Value = 642
$ws.Range("K6").Value = 635
$ws.Range("L6").Value = 637
$ws.Range("M6").Value = 593
$ws.Range("H7").Value = 535
$ws.Range("I7").Value = 536
$ws.Range("J7").Value = 577
$ws.Range("K7").Value = 499
$ws.Range("B8").Value = 532
$ws.Range("C8").Value = 530
$ws.Range("D8").Value = 530
$ws.Range("E8").Value = 534
$ws.Range("F8").Value = 542
$ws.Range("G8").Value = 529
$ws.Range("H8").Value = 754
$ws.Range("I8").Value = 731
$ws.Range("J8").Value = 684
$ws.Range("K8").Value = 649
$ws.Range("L8").Value = 631
$ws.Range("M8").Value = 633
$ws.Range("B9").Value = 415
$ws.Range("C9").Value = 513
$ws.Range("D9").Value = 512
$ws.Range("E9").Value = 519
$ws.Range("F9").Value = 521
$ws.Range("G9").Value = 499
$ws.Range("I9").Value = 734
$ws.Range("J9").Value = 711
$ws.Range("K9").Value = 671
$ws.Range("L9").Value = 664
$ws.Range("M9").Value = 652
$ws.Range("C10").Value = 441
$ws.Range("D10").Value = 422
$ws.Range("E10").Value = 427
$ws.Range("F10").Value = 448
$ws.Range("H10").Value = 631
$ws.Range("I10").Value = 602
$ws.Range("J10").Value = 572
$ws.Range("K10").Value = 576
$ws.Range("L10").Value = 542
$ws.Range("M10").Value = 532
$ws.Range("C11").Value = 331
$ws.Range("E11").Value = 381
$ws.Range("F11").Value = 366
$ws.Range("G11").Value = 344
$ws.Range("H11").Value = 598
$ws.Range("I11").Value = 591
$ws.Range("J11").Value = 577
$ws.Range("K11").Value = 567
$ws.Range("L11").Value = 559
$ws.Range("M11").Value = 551
$ws.Range("I12").Value = 544
$ws.Range("J12").Value = 499
$ws.Range("K12").Value = 528
$ws.Range("L12").Value = 542
$ws.Range("M12").Value = 462
$ws.Range("H13").Value = 486
$ws.Range("I13").Value = 480
$ws.Range("J13").Value = 437
$ws.Range("K13").Value = 449
$ws.Range("L13").Value = 379
$ws.Range("C14").Value = 448
$ws.Range("D14").Value = 426
$ws.Range("E14").Value = 439
$ws.Range("F14").Value = 440
$ws.Range("G14").Value = 411
$ws.Range("H14").Value = 565
$ws.Range("I14").Value = 550
$ws.Range("J14").Value = 557
$ws.Range("K14").Value = 546
$ws.Range("L14").Value = 576
$ws.Range("M14").Value = 542
$ws.Range("C15").Value = 333
$ws.Range("D15").Value = 359
$ws.Range("E15").Value = 362
$ws.Range("F15").Value = 366
$ws.Range("G15").Value = 370
$ws.Range("H15").Value = 558
$ws.Range("I15").Value = 556
$ws.Range("J15").Value = 515
$ws.Range("K15").Value = 540
$ws.Range("M15").Value = 550
$ws.Range("D16").Value = 285
$ws.Range("E16").Value = 295
$ws.Range("F16").Value = 342
$ws.Range("G16").Value = 247
$ws.Range("H16").Value = 461
$ws.Range("I16").Value = 447
$ws.Range("J16").Value = 455
$ws.Range("K16").Value = 435
$ws.Range("L16").Value = 438
$ws.Range("M16").Value = 445
$ws.Range("D17").Value = 0
$ws.Range("H17").Value = 446
$ws.Range("I17").Value = 435
$ws.Range("J17").Value = 420
$ws.Range("K17").Value = 407
$ws.Range("L17").Value = 412
$ws.Range("M17").Value = 424
$ws.Range("H18").Value = 370
$ws.Range("I18").Value = 378
$ws.Range("J18").Value = 394
$ws.Range("K18").Value = 370
$ws.Range("L18").Value = 382
$ws.Range("M18").Value = 382
$ws.Range("H19").Value = 257
$ws.Range("I19").Value = 252
$ws.Range("J19").Value = 258
$ws.Range("K19").Value = 268
$ws.Range("L19").Value = 301
$ws.Range("M19").Value = 325
$ws.Range("C20").Value = 315
$ws.Range("D20").Value = 449
$ws.Range("E20").Value = 422
$ws.Range("F20").Value = 477
$ws.Range("G20").Value = 462
$ws.Range("H20").Value = 748
$ws.Range("I20").Value = 758
$ws.Range("J20").Value = 774
$ws.Range("K20").Value = 798
$ws.Range("L20").Value = 825
$ws.Range("M20").Value = 786
$ws.Range("C21").Value = 244
$ws.Range("D21").Value = 408
$ws.Range("E21").Value = 450
$ws.Range("F21").Value = 542
$ws.Range("G21").Value = 609
$ws.Range("H21").Value = 932
$ws.Range("I21").Value = 867
$ws.Range("J21").Value = 864
$ws.Range("K21").Value = 834
$ws.Range("L21").Value = 928
$ws.Range("M21").Value = 958
$ws.Range("B22").Value = 322
$ws.Range("C22").Value = 217
$ws.Range("D22").Value = 352
$ws.Range("E22").Value = 416
$ws.Range("F22").Value = 461
$ws.Range("G22").Value = 464
$ws.Range("H22").Value = 753
$ws.Range("I22").Value = 721
$ws.Range("J22").Value = 602
$ws.Range("K22").Value = 508
$ws.Range("L22").Value = 496
$ws.Range("M22").Value = 806
$ws.Range("C23").Value = 224
$ws.Range("D23").Value = 225
$ws.Range("E23").Value = 229
$ws.Range("F23").Value = 240
$ws.Range("H23").Value = 674
$ws.Range("I23").Value = 603
$ws.Range("J23").Value = 694
$ws.Range("K23").Value = 493
$ws.Range("L23").Value = 563
$ws.Range("M23").Value = 982
